$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 497..499 (everything from old row 497 onward shifts down by 3,
# dimension grows from A1:T577 to A1:T580). The new rows carry the newest week's
# data for this product/terminal, inserted above the existing (now-shifted) rows.
$ws.Rows("497:499").Insert()

# Row 497 - Pintón
$ws.Range("A497").Value = 8
$ws.Range("B497").Value = "Terminal La Palmera de La Serena"
$ws.Range("C497").Value = "Coquimbo"
$ws.Range("D497").Value = 44637
$ws.Range("E497").Value = 4
$ws.Range("F497").Value = "Fruta"
$ws.Range("G497").Value = 100108
$ws.Range("H497").Value = "Tropicales y subtropicales"
$ws.Range("I497").Value = 100108006
$ws.Range("J497").Value = "Plátano"
$ws.Range("K497").Value = "Sin especificar"
$ws.Range("L497").Value = "Pintón"
$ws.Range("M497").Value = 80
$ws.Range("N497").Value = 18000
$ws.Range("O497").Value = 18000
$ws.Range("P497").Value = 18000
$ws.Range("Q497").Value = "$/caja 20 kilos"
$ws.Range("R497").Value = "Ecuador"
$ws.Range("S497").Value = 900
$ws.Range("T497").Value = 20

# Row 498 - Primera Maduro
$ws.Range("A498").Value = 8
$ws.Range("B498").Value = "Terminal La Palmera de La Serena"
$ws.Range("C498").Value = "Coquimbo"
$ws.Range("D498").Value = 44637
$ws.Range("E498").Value = 4
$ws.Range("F498").Value = "Fruta"
$ws.Range("G498").Value = 100108
$ws.Range("H498").Value = "Tropicales y subtropicales"
$ws.Range("I498").Value = 100108006
$ws.Range("J498").Value = "Plátano"
$ws.Range("K498").Value = "Sin especificar"
$ws.Range("L498").Value = "Primera Maduro"
$ws.Range("M498").Value = 120
$ws.Range("N498").Value = 20000
$ws.Range("O498").Value = 20000
$ws.Range("P498").Value = 20000
$ws.Range("Q498").Value = "$/caja 20 kilos"
$ws.Range("R498").Value = "Ecuador"
$ws.Range("S498").Value = 1000
$ws.Range("T498").Value = 20

# Row 499 - Primera Pintón
$ws.Range("A499").Value = 8
$ws.Range("B499").Value = "Terminal La Palmera de La Serena"
$ws.Range("C499").Value = "Coquimbo"
$ws.Range("D499").Value = 44637
$ws.Range("E499").Value = 4
$ws.Range("F499").Value = "Fruta"
$ws.Range("G499").Value = 100108
$ws.Range("H499").Value = "Tropicales y subtropicales"
$ws.Range("I499").Value = 100108006
$ws.Range("J499").Value = "Plátano"
$ws.Range("K499").Value = "Sin especificar"
$ws.Range("L499").Value = "Primera Pintón"
$ws.Range("M499").Value = 120
$ws.Range("N499").Value = 21000
$ws.Range("O499").Value = 21000
$ws.Range("P499").Value = 21000
$ws.Range("Q499").Value = "$/caja 20 kilos"
$ws.Range("R499").Value = "Ecuador"
$ws.Range("S499").Value = 1050
$ws.Range("T499").Value = 20
